# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the refreshed data snapshot (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - F-column row -> new value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 5618
$ws1.Range("F4").Value  = 7664
$ws1.Range("F6").Value  = 77
$ws1.Range("F8").Value  = 613
$ws1.Range("F11").Value = 4435
$ws1.Range("F14").Value = 119
$ws1.Range("F15").Value = 2996
$ws1.Range("F17").Value = 573
$ws1.Range("F19").Value = 548
$ws1.Range("F20").Value = 476
$ws1.Range("F21").Value = 482
$ws1.Range("F24").Value = 1726
$ws1.Range("F27").Value = 1454
$ws1.Range("F37").Value = 3139
$ws1.Range("F38").Value = 719
$ws1.Range("F39").Value = 47
$ws1.Range("F40").Value = 161
$ws1.Range("F41").Value = 48
$ws1.Range("F42").Value = 1002

# Sheet "全部类型" (all types) - same underlying rows, shifted by the extra
# "本地生活" row that appears at row 22 on this combined sheet.
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 5618
$ws4.Range("F4").Value  = 7664
$ws4.Range("F6").Value  = 77
$ws4.Range("F8").Value  = 613
$ws4.Range("F11").Value = 4435
$ws4.Range("F14").Value = 119
$ws4.Range("F15").Value = 2996
$ws4.Range("F17").Value = 573
$ws4.Range("F19").Value = 548
$ws4.Range("F20").Value = 476
$ws4.Range("F21").Value = 482
$ws4.Range("F25").Value = 1726
$ws4.Range("F28").Value = 1454
$ws4.Range("F38").Value = 3139
$ws4.Range("F40").Value = 719
$ws4.Range("F41").Value = 47
$ws4.Range("F42").Value = 161
$ws4.Range("F43").Value = 48
$ws4.Range("F44").Value = 1002

$wb.Save()
